# Insert a new "Vlak" worksheet right after "Objecten", matching the
# Put-sheet layout/template but describing the Vlak (level) attribute.
$wb = $excel.ActiveWorkbook

$objecten = $wb.Worksheets.Item("Objecten")
$vlak = $wb.Worksheets.Add($null, $objecten)
$vlak.Name = "Vlak"

$vlak.Range("A1").Value = "Attribute"
$vlak.Range("B1").Value = "Kolommen"
$vlak.Range("A2").Value = "putnr"
$vlak.Range("B2").Value = '["PUT", "PUTNO"]'
$vlak.Range("A3").Value = "vlaknr"
$vlak.Range("B3").Value = '["VLAK", "VLAKNO"]'

$vlak.Columns.Item(2).ColumnWidth = 19.83

# Update selection on "Put" sheet (it moved down a slot but keeps its content)
$put = $wb.Worksheets.Item("Put")
$put.Activate()
$put.Range("A2:B2").Select()

# Make "Vlak" the active/visible tab with its own zoom + selection
$vlak.Activate()
$excel.ActiveWindow.Zoom = 140
$vlak.Range("D10").Select()
